$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row additions: P1 = 14, Q1 = 15
$ws.Cells.Item(1, 16).Value = 14
$ws.Cells.Item(1, 17).Value = 15

# Copy formatting (border/bold/alignment) from O1 to the new header cells
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Updated data values for columns B:G, I:K (rows 2-25); column H is cleared (deleted)
$colIndexes = @(2,3,4,5,6,7,9,10,11)
$data = @(
    @(24.94308144648824, 21.58519479409142, 3.172920026738549, 30.01148016419559, 22.43282096852754, 31.28937338776164, 3.955968496043898, 9.575560831914753, 12.61737568504521),
    @(23.28835297949403, 20.17920820822495, 3.169207893864802, 27.79961142433553, 21.38414764847752, 29.56294045941251, 3.677163409813299, 9.40843685713501, 12.67875287897997),
    @(22.21366093938064, 19.27554437604578, 3.167561397034676, 26.3807008810615, 20.73882656332893, 28.49083949975601, 3.502426779641838, 9.313597194351102, 12.72990372569956),
    @(21.76058962036348, 18.90865478611785, 3.168777135723311, 25.78665342752303, 20.45706410749041, 28.01738293302167, 3.430361273174948, 9.270207986549163, 12.7407609722592),
    @(21.68412198155679, 18.86128961335256, 3.171121974242187, 25.68690209249982, 20.38738756989021, 27.89705099453713, 3.418606895448627, 9.255011420978049, 12.72675892881022),
    @(22.20671007374737, 19.30845431302579, 3.173453669286967, 26.37233537825206, 20.67345545494991, 28.37302875873894, 3.50189490655549, 9.290980109979827, 12.68626698839713),
    @(24.38352783020316, 21.15501596777912, 3.179537711700251, 29.26059751751236, 21.99467563882916, 30.55896430617418, 3.860647343248495, 9.487726149850612, 12.57732451251094),
    @(28.19290278116164, 24.39980507880212, 3.185797163727982, 34.46335409541976, 24.66647323454185, 34.91750631635735, 4.539124298408021, 9.976749768835369, 12.5634867217338),
    @(30.7177510810272, 26.6140271731334, 3.197344917956122, 38.03608087990279, 26.58286553855087, 37.99882132270059, 5.026772231314881, 10.3696960400201, 12.62532195421664),
    @(31.80761154670709, 27.62775245257356, 3.21451956322756, 39.61459479175513, 27.37368470365023, 39.3061726052534, 5.247907975071698, 10.52734875740375, 12.60751724561919),
    @(32.21286851605198, 27.98169516408214, 3.215427208973541, 40.20704301757286, 27.71480620257323, 39.87450230357955, 5.331295260664352, 10.60584638104839, 12.64273249563317),
    @(32.12607682766924, 27.90088309999386, 3.214053123734724, 40.07979267575375, 27.64961050230452, 39.76644550921134, 5.313245367779563, 10.59226526260631, 12.6425513038555),
    @(31.84116902391085, 27.65491356591168, 3.214094326524951, 39.66348909204493, 27.40529139064379, 39.35908239278616, 5.254724857071578, 10.53522437234301, 12.61356013115944),
    @(31.66534785860943, 27.51320783584488, 3.216439869309077, 39.40756987646571, 27.23906976576681, 39.08066994623039, 5.219122320464706, 10.49374360403099, 12.58135243535914),
    @(30.6431933111296, 26.6301060325048, 3.215413961501683, 37.93047256997998, 26.38895993578022, 37.66865197801809, 5.01367132217232, 10.30196400177566, 12.50022008070711),
    @(30.00152505067592, 26.07887770587025, 3.214928806094902, 37.0131459597586, 25.86616046298857, 36.82628832822994, 4.8875387393154, 10.18717068166766, 12.45767192010447),
    @(29.62769909805493, 25.7323090857154, 3.208908008292842, 36.48162672126502, 25.60929671393956, 36.41776614663112, 4.813884233671957, 10.13990972763268, 12.47421354659291),
    @(29.49967896165456, 25.63433557569831, 3.211344359569056, 36.30053548062338, 25.48813177869704, 36.21946405629559, 4.789676767472248, 10.11041558183678, 12.45055694837129),
    @(30.07044353666289, 26.13590917650443, 3.214518425202436, 37.11127921146612, 25.92530204158497, 36.92213596769967, 4.900887890944042, 10.20066739461202, 12.4649524419607),
    @(31.92434408545504, 27.75441369882138, 3.220656800230934, 39.78535661843782, 27.43063738556625, 39.39877321697031, 5.272616750242205, 10.53283214101443, 12.5793778067166),
    @(33.08976857190756, 28.74765756703172, 3.216815210856497, 41.50055003039719, 28.54800207785669, 41.12836651668138, 5.514690681442437, 10.78383497779948, 12.73336906892805),
    @(32.47289076644926, 28.1878794981463, 3.210881877772938, 40.58854524124547, 27.97689313152366, 40.30232215915805, 5.38486399581974, 10.67184370464864, 12.69998794225866),
    @(30.04071046394963, 26.0558428187124, 3.202447085851767, 37.06792092265584, 25.99132841598003, 37.04097832223374, 4.89331828954009, 10.23188189785393, 12.54306077908021),
    @(27.20954548668338, 23.61223895018857, 3.194973080784538, 33.10133983697861, 23.85010757766425, 33.58039591958165, 4.358482535645046, 9.799433413823545, 12.47265843523525)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $colIndexes.Count; $j++) {
        $ws.Cells.Item($row, $colIndexes[$j]).Value = $rowVals[$j]
    }
}

# Remove column H data (rows 2-25) entirely
$ws.Range("H2:H25").ClearContents()

# New columns P and Q (rows 2-25) with zero values
for ($row = 2; $row -le 25; $row++) {
    $ws.Cells.Item($row, 16).Value = 0
    $ws.Cells.Item($row, 17).Value = 0
}

Write-Host "Edit complete"